# "added name to report"
#
# The author's name placeholder on the title page ("Alexander Simes
# ???????") is completed with the real student number (4415299),
# matching the other two authors' "Name NUMBER" pattern already present
# in the document. Word also drops the "_GoBack" bookmark at the point
# of the last edit, which we reproduce here too.

$d = $word.ActiveDocument

# Step 1: collapse "Alexander Simes ???????" down to "Alexander Simes "
# (i.e. drop the "???????" placeholder, keeping a single trailing space).
$rng = $d.Content
$rng.Find.Execute("Alexander Simes ???????", $false, $false, $false, $false, `
                   $false, $true, 1, $false, "", 0)
$rng.Text = "Alexander Simes "

# Step 2: locate that text again, move to just after it (right before the
# following page break) and type in the student number there.
$rng2 = $d.Content
$rng2.Find.Execute("Alexander Simes ", $false, $false, $false, $false, `
                    $false, $true, 1, $false, "", 0)
$rng2.Collapse(0)

# Word stamps the last edit location with the "_GoBack" bookmark.
$d.Bookmarks.Add("_GoBack", $rng2)
$numRng = $d.Bookmarks.Item("_GoBack").Range
$numRng.InsertBefore("4415299")
